# Update odds values on Sheet1 (rows 3, 6, 7) to reflect the latest
# FlashScore scrape, per the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Guarani vs CRB) ---
$ws.Cells.Item(3, 7).Value = 2.1      # G3  Odd_H_FT
$ws.Cells.Item(3, 9).Value = 4.1      # I3  Odd_A_FT
$ws.Cells.Item(3, 10).Value = 3       # J3  Odd_H_HT
$ws.Cells.Item(3, 12).Value = 5       # L3  Odd_A_HT
$ws.Cells.Item(3, 21).Value = 2.5     # U3  Odd_BTTS_Yes
$ws.Cells.Item(3, 22).Value = 1.5     # V3  Odd_BTTS_No
$ws.Cells.Item(3, 23).Value = 5       # W3  Odd_CS_1-0
$ws.Cells.Item(3, 25).Value = 10      # Y3  Odd_CS_2-1
$ws.Cells.Item(3, 26).Value = 19      # Z3  Odd_CS_3-0
$ws.Cells.Item(3, 32).Value = 101     # AF3 Odd_CS_3-3
$ws.Cells.Item(3, 33).Value = 8       # AG3 Odd_CS_0-1
$ws.Cells.Item(3, 34).Value = 19      # AH3 Odd_CS_0-2
$ws.Cells.Item(3, 35).Value = 17      # AI3 Odd_CS_1-2
$ws.Cells.Item(3, 50).Value = 26      # AX3 Odd_CS_0-2_HT
$ws.Cells.Item(3, 52).Value = 101     # AZ3 Odd_CS_0-3_HT

# --- Row 6 ---
$ws.Cells.Item(6, 7).Value = 2.1      # G6  Odd_H_FT
$ws.Cells.Item(6, 8).Value = 3.2      # H6  Odd_D_FT
$ws.Cells.Item(6, 9).Value = 3.7      # I6  Odd_A_FT
$ws.Cells.Item(6, 10).Value = 2.88    # J6  Odd_H_HT
$ws.Cells.Item(6, 12).Value = 4.33    # L6  Odd_A_HT
$ws.Cells.Item(6, 17).Value = 2.35    # Q6  Odd_Over25_FT
$ws.Cells.Item(6, 18).Value = 1.57    # R6  Odd_Under25_FT
$ws.Cells.Item(6, 21).Value = 2.1     # U6  Odd_BTTS_Yes
$ws.Cells.Item(6, 22).Value = 1.67    # V6  Odd_BTTS_No
$ws.Cells.Item(6, 23).Value = 6       # W6  Odd_CS_1-0
$ws.Cells.Item(6, 24).Value = 9       # X6  Odd_CS_2-0
$ws.Cells.Item(6, 25).Value = 9.5     # Y6  Odd_CS_2-1
$ws.Cells.Item(6, 26).Value = 19      # Z6  Odd_CS_3-0
$ws.Cells.Item(6, 27).Value = 19      # AA6 Odd_CS_3-1
$ws.Cells.Item(6, 28).Value = 34      # AB6 Odd_CS_3-2
$ws.Cells.Item(6, 30).Value = 6       # AD6 Odd_CS_1-1
$ws.Cells.Item(6, 31).Value = 19      # AE6 Odd_CS_2-2
$ws.Cells.Item(6, 32).Value = 67      # AF6 Odd_CS_3-3
$ws.Cells.Item(6, 33).Value = 8.5     # AG6 Odd_CS_0-1
$ws.Cells.Item(6, 34).Value = 17      # AH6 Odd_CS_0-2
$ws.Cells.Item(6, 35).Value = 13      # AI6 Odd_CS_1-2
$ws.Cells.Item(6, 36).Value = 41      # AJ6 Odd_CS_0-3
$ws.Cells.Item(6, 37).Value = 34      # AK6 Odd_CS_1-3
$ws.Cells.Item(6, 38).Value = 41      # AL6 Odd_CS_2-3
$ws.Cells.Item(6, 40).Value = 4       # AN6 Odd_CS_1-0_HT
$ws.Cells.Item(6, 41).Value = 12      # AO6 Odd_CS_2-0_HT
$ws.Cells.Item(6, 43).Value = 41      # AQ6 Odd_CS_3-0_HT
$ws.Cells.Item(6, 45).Value = 201     # AS6 Odd_CS_3-2_HT
$ws.Cells.Item(6, 47).Value = 9       # AU6 Odd_CS_1-1_HT
$ws.Cells.Item(6, 48).Value = 67      # AV6 Odd_CS_2-2_HT
$ws.Cells.Item(6, 49).Value = 5.5     # AW6 Odd_CS_0-1_HT
$ws.Cells.Item(6, 50).Value = 21      # AX6 Odd_CS_0-2_HT
$ws.Cells.Item(6, 51).Value = 34      # AY6 Odd_CS_1-2_HT
$ws.Cells.Item(6, 52).Value = 81      # AZ6 Odd_CS_0-3_HT
$ws.Cells.Item(6, 53).Value = 101     # BA6 Odd_CS_1-3_HT
$ws.Cells.Item(6, 54).Value = 301     # BB6 Odd_CS_2-3_HT

# --- Row 7 ---
$ws.Cells.Item(7, 15).Value = 1.29    # O7 Odd_Over15_FT
$ws.Cells.Item(7, 16).Value = 3.5     # P7 Odd_Under15_FT
$ws.Cells.Item(7, 17).Value = 1.98    # Q7 Odd_Over25_FT
$ws.Cells.Item(7, 18).Value = 1.88    # R7 Odd_Under25_FT
